# Applies the "Append: 2025-10-12 12:31 JST" change to the "ランサーズ" sheet.
# New scrape results are inserted above the previously-seen rows (rows shift down),
# all "取得日時" timestamps are refreshed, and column H is widened.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Full target content for data rows 2-10 (row 1 header is unchanged).
$rowsData = @(
  @{ Row=2; A="2025-10-12 12:31:05"; B="急募 PR Zoom/Meet×TLDV×ChatGPT×Notion×Slack 議事録ワークフロー構築依頼"; C="システム開発"; D="50,000 円 ~ 100,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5410688"; G=323; H="🔥GPT,ChatGPT" },
  @{ Row=3; A="2025-10-12 12:31:05"; B="【急募】Inkscape「Hershey Text」用svgフォント変換ツール開発(python)"; C="システム開発"; D="10,000 円 ~ 20,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5411941"; G=315; H="🔥Python ◆ツール,開発" },
  @{ Row=4; A="2025-10-12 12:31:05"; B="【個人PoC案件】ChatGPT(web)とAzure連携の仕組み構築"; C="システム開発"; D="10,000 円 ~ 20,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5411897"; G=315; H="🔥GPT,ChatGPT" },
  @{ Row=5; A="2025-10-12 12:31:05"; B="【急募】予定管理のWebシステム開発をお手伝いください!"; C="システム開発"; D="100,000 円 ~ 200,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5411923"; G=153; H="◆開発,システム開発 ◇管理" },
  @{ Row=6; A="2025-10-12 12:31:05"; B="【自動売買】Excelと楽天RSSを活用したシステム開発依頼"; C="システム開発"; D="5,000 円 ~ 10,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5411684"; G=110; H="◆開発,システム開発" },
  @{ Row=7; A="2025-10-12 12:31:05"; B="Laravelでのバックエンド開発:管理画面機能やDB管理・ポイント機能などの開発【フルリモート】"; C="システム開発"; D="200,000 円 ~ 300,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5411736"; G=93; H="◆開発 ◇管理" },
  @{ Row=8; A="2025-10-12 12:31:05"; B="スプレッドシートをもとにした顧客・売上管理アプリのグライド化(Glide/無料版)"; C="システム開発"; D="5,000 円 ~ 10,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5411871"; G=55; H="◇アプリ" },
  @{ Row=9; A="2025-10-12 12:31:05"; B="【急募】教育系のWEBサイトの作成"; C="システム開発"; D="20,000 円 ~ 50,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5411679"; G=33; H="◇サイト" },
  @{ Row=10; A="2025-10-12 12:31:05"; B="【急募】微生物の画像判定を行う専門家を探しています!"; C="システム開発"; D="50,000 円 ~ 100,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5411887"; G=18; H=$null }
)

# Clear any hyperlinks currently on the sheet; they will be rebuilt below for the new row layout.
$ws.Range("F2").Hyperlinks.Delete()

foreach ($item in $rowsData) {
  $r = $item.Row
  $ws.Cells.Item($r, 1).Value = $item.A
  $ws.Cells.Item($r, 2).Value = $item.B
  $ws.Cells.Item($r, 3).Value = $item.C
  $ws.Cells.Item($r, 4).Value = $item.D
  $ws.Cells.Item($r, 5).Value = $item.E
  $ws.Cells.Item($r, 6).Value = $item.F
  $ws.Cells.Item($r, 7).Value = $item.G
  if ($item.H -ne $null) {
    $ws.Cells.Item($r, 8).Value = $item.H
  }
  $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $item.F)
  $ws.Cells.Item($r, 6).Style = "Hyperlink"
}

# Column H ("スキル概要") grows from width 14 to width 17 to fit the new entries.
$ws.Columns.Item(8).ColumnWidth = 16.17

Write-Output "applied"
